# Daily attendance processing - 2026-01-02 21:54:49
# Reorders the "Recorded By" (column G) contributor lists for specific rows
# on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G2"   = "System, backup@backdoor.com, system"
    "G3"   = "System, dnasr281@gmail.com"
    "G4"   = "System, backup@backdoor.com"
    "G5"   = "System, backup@backdoor.com"
    "G6"   = "System, dnasr281@gmail.com"
    "G7"   = "System, admin@admin.com"
    "G8"   = "System, backup@backdoor.com"
    "G28"  = "System, backup@backdoor.com, system"
    "G29"  = "System, dnasr281@gmail.com"
    "G30"  = "System, backup@backdoor.com"
    "G31"  = "System, backup@backdoor.com"
    "G32"  = "System, dnasr281@gmail.com"
    "G33"  = "System, admin@admin.com"
    "G34"  = "System, backup@backdoor.com"
    "G54"  = "System, backup@backdoor.com, system"
    "G55"  = "System, dnasr281@gmail.com"
    "G56"  = "System, backup@backdoor.com"
    "G57"  = "System, backup@backdoor.com"
    "G58"  = "System, dnasr281@gmail.com"
    "G59"  = "System, admin@admin.com"
    "G60"  = "System, backup@backdoor.com"
    "G80"  = "System, backup@backdoor.com"
    "G81"  = "System, backup@backdoor.com"
    "G82"  = "System, backup@backdoor.com"
    "G87"  = "dnasr281@gmail.com, admin@admin.com"
    "G106" = "System, backup@backdoor.com"
    "G107" = "System, backup@backdoor.com"
    "G108" = "System, backup@backdoor.com"
    "G113" = "dnasr281@gmail.com, admin@admin.com"
    "G132" = "System, backup@backdoor.com"
    "G133" = "System, backup@backdoor.com"
    "G134" = "System, backup@backdoor.com"
    "G139" = "dnasr281@gmail.com, admin@admin.com"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
